# unitTest_numberCommand.xlsx - add new "web" command `deselect(locator,text)`
#
# The hidden "#system" sheet keeps, in column U, an alphabetically sorted
# legend of all `web` commands (used to populate the `web` named range /
# dropdown validation). We need to insert a new entry -
# "deselect(locator,text)" - right before the existing
# "deselectMulti(locator,array)" entry (row 53), pushing every entry from
# old U53 down through old U116 to U54..U117. Only column U moves; every
# other column (e.g. E, which holds an unrelated "desktop" legend) must
# stay exactly where it is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$firstRow = 53
$lastRow = 116
$col = 21   # column U

# Shift existing values down by one row, working from the bottom up so we
# don't clobber a value before it has been copied.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $val = $ws.Cells.Item($r, $col).Value()
    $ws.Cells.Item($r + 1, $col).Value = $val
}

# Insert the new command text in the now-vacated row.
$ws.Cells.Item($firstRow, $col).Value = "deselect(locator,text)"

# Extend the "web" named range by one row to cover the new entry.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "web") {
        $n.RefersTo = "='#system'!`$U`$2:`$U`$117"
    }
}
